$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CICIDS20180")
$ws.Range("B2").Value = 0.9919494986534119
$ws.Range("C2").Value = 0.9791200160980225
$ws.Range("D2").Value = 10

$ws = $wb.Worksheets.Item("CICIDS2018_0")
$ws.Range("B13").Value = 0.9919494986534119
$ws.Range("C13").Value = 0.9791200160980225
$ws.Range("D13").Value = 10
$ws.Range("B14").Value = 0.9919494986534119
$ws.Range("C14").Value = 0.9791200160980225
$ws.Range("D14").Value = 10
$ws.Range("B15").Value = 0.9919494986534119
$ws.Range("C15").Value = 0.9791200160980225
$ws.Range("D15").Value = 10

$ws = $wb.Worksheets.Item("CICIDS20181")
$ws.Range("B2").Value = 0.9919494986534119
$ws.Range("C2").Value = 0.9791200160980225
$ws.Range("D2").Value = 10
$ws.Range("B3").Value = 0.01773619651794434
$ws.Range("C3").Value = 0.9948400259017944
$ws.Range("D3").Value = 10

$ws = $wb.Worksheets.Item("CICIDS2018_1")
$ws.Range("B13").Value = 0.01773619651794434
$ws.Range("C13").Value = 0.9791200160980225
$ws.Range("D13").Value = 10
$ws.Range("B14").Value = 0.9919494986534119
$ws.Range("C14").Value = 0.9948400259017944
$ws.Range("D14").Value = 10
$ws.Range("B15").Value = 0.5048428475856781
$ws.Range("C15").Value = 0.9869800209999084
$ws.Range("D15").Value = 10
$ws.Range("B16").Value = 0.6888728322621279
$ws.Range("C16").Value = 0.01111572553256617

$ws = $wb.Worksheets.Item("CICIDS20182")
$ws.Range("B2").Value = 0.9919494986534119
$ws.Range("C2").Value = 0.9791200160980225
$ws.Range("D2").Value = 10
$ws.Range("B3").Value = 0.01773619651794434
$ws.Range("C3").Value = 0.9948400259017944
$ws.Range("D3").Value = 10
$ws.Range("B4").Value = 0.0266546867787838
$ws.Range("C4").Value = 0.9916239380836487
$ws.Range("D4").Value = 10

$ws = $wb.Worksheets.Item("CICIDS2018_2")
$ws.Range("B13").Value = 0.01773619651794434
$ws.Range("C13").Value = 0.9791200160980225
$ws.Range("D13").Value = 10
$ws.Range("B14").Value = 0.9919494986534119
$ws.Range("C14").Value = 0.9948400259017944
$ws.Range("D14").Value = 10
$ws.Range("B15").Value = 0.34544679398338
$ws.Range("C15").Value = 0.9885279933611552
$ws.Range("D15").Value = 10
$ws.Range("B16").Value = 0.5599055234826524
$ws.Range("C16").Value = 0.008304717475701578
